# Apply data updates to Sheet1 as described in the commit diff
# ("Updated symbol list on Tue Dec 27 05:55:11 UTC 2022 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.34"
$ws.Range("D3").Value = "'23.08"
$ws.Range("D4").Value = "'5.408"
$ws.Range("D7").Value = "'6.493"
$ws.Range("D8").Value = "'0.8092"
$ws.Range("D9").Value = "'0.9249"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01112"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1427"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07383"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03274"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03074"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09353"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.868"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001571"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04715"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "'0.005848"
$ws.Range("D21").Value = "'0.004878"
$ws.Range("D22").Value = "'0.00006797"
$ws.Range("D23").Value = "'3.574"
$ws.Range("D24").Value = "'2.178"
$ws.Range("D25").Value = "'0.3236"
$ws.Range("D40").Value = "'0.03965"
$ws.Range("D41").Value = "'0.006367"
$ws.Range("D42").Value = "'0.004299"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("D43").Value = "'0.1076"
$ws.Range("D44").Value = "'0.008452"
$ws.Range("D45").Value = "'0.00005091"
$ws.Range("D47").Value = "'0.6698"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Value = "'0.002449"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D50").Value = "'0.0001999"
